# Add renda familiar em valor discreto
# Delete the "Q4" row (row 4) entirely - the subsequent rows shift up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:M4").EntireRow.Delete() | Out-Null

# Update remaining rows with recalculated statistics
# Row 2 - Q1
$ws.Range("C2").Value = 22.91
$ws.Range("E2").Value = 18
$ws.Range("I2").Value = 3.91
$ws.Range("J2").Value = 15.31
$ws.Range("K2").Value = 0.17
$ws.Range("L2").Value = "{'Q1': np.float64(19.0), 'Q2': np.float64(23.0), 'Q3': np.float64(26.0), 'Q4': np.float64(29.0)}"
$ws.Range("M2").Value = 7

# Row 3 - Q3
$ws.Range("C3").Value = 6.2
$ws.Range("D3").Value = 6
$ws.Range("I3").Value = 3.51
$ws.Range("J3").Value = 12.29
$ws.Range("K3").Value = 0.57
$ws.Range("L3").Value = "{'Q1': np.float64(3.0), 'Q2': np.float64(6.0), 'Q3': np.float64(9.0), 'Q4': np.float64(12.0)}"
$ws.Range("M3").Value = 6

# Row 4 (was row 5) - Q7
$ws.Range("C4").Value = 4.33
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 2.58
$ws.Range("J4").Value = 6.67
$ws.Range("K4").Value = 0.6
$ws.Range("L4").Value = "{'Q1': np.float64(2.0), 'Q2': np.float64(4.0), 'Q3': np.float64(7.0), 'Q4': np.float64(8.0)}"
$ws.Range("M4").Value = 5

# Row 5 (was row 6) - Q9
$ws.Range("C5").Value = 13.23
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 20
$ws.Range("I5").Value = 6.45
$ws.Range("J5").Value = 41.58
$ws.Range("K5").Value = 0.49
$ws.Range("L5").Value = "{'Q1': np.float64(7.0), 'Q2': np.float64(13.0), 'Q3': np.float64(19.0), 'Q4': np.float64(25.0)}"
$ws.Range("M5").Value = 12

# Row 6 (was row 7) - Q11
$ws.Range("E6").Value = 3
$ws.Range("J6").Value = 6.71

# Row 7 (was row 8) - Q13
$ws.Range("C7").Value = 5.65
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 10
$ws.Range("I7").Value = 2.97
$ws.Range("J7").Value = 8.789999999999999
$ws.Range("K7").Value = 0.53
$ws.Range("L7").Value = "{'Q1': np.float64(3.0), 'Q2': np.float64(6.0), 'Q3': np.float64(8.0), 'Q4': np.float64(10.0)}"
